$p = $ppt.ActivePresentation
Write-Output "Slides count: $($p.Slides.Count)"
